$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ls_hvcb_dp")

$ws.Range("F54").Value = "ssun_132"
$ws.Range("F55").Value = "ssun_132"

$ws.Range("F75").Value = "kbek_132"
$ws.Range("F76").Value = "kbek_132"

$ws.Range("F97").Value = "nsci_132"
$ws.Range("F98").Value = "nsci_132"

$ws.Range("F115").Value = "xnss_132"
$ws.Range("F116").Value = "xnss_132"

$ws.Range("F117").Value = "lkjg_132"
$ws.Range("F118").Value = "lkjg_132"

$ws.Range("F181").Value = "ksni_132"
$ws.Range("F182").Value = "ksni_132"

$ws.Range("F189").Value = "utmj_132"
$ws.Range("F190").Value = "utmj_132"

$ws.Range("F201").Value = "kcmt_132"
$ws.Range("F202").Value = "kcmt_132"

$ws.Range("F203").Value = "tpau_132"

$ws.Range("F269").Value = "pltg_230"
$ws.Range("F270").Value = "pltg_230"

[void]$ws.Range("F271").Select()
